$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header tweak: P1 "CSI Camera" -> "Camera" ---
$ws.Range("P1").Value = "Camera"

# --- Move a few stray values out of the to-be-deleted X/Y columns ---
# Row 11: camera info that belonged in P (CSI Camera / Camera) column.
$ws.Range("P11").Value = "HSB camera via QSFP slot USB camera"

# Rows 12-13: same camera description duplicated.
$cameraText = "Up to 20 cameras via HSB Up to 6 cameras through 16x lanes MIPI CSI-2 Up to 32 cameras using Virtual Channels C-PHY 2.1 (10.25 Gbps) D-PHY 2.1 (40 Gbps)"
$ws.Range("P12").Value = $cameraText
$ws.Range("P13").Value = $cameraText

# Rows 14-18: DL Accelerator value that belonged in H column.
$ws.Range("H14").Value = "2x NVDLA"
$ws.Range("H15").Value = "2x NVDLA"
$ws.Range("H16").Value = "2x NVDLA"
$ws.Range("H17").Value = "2x NVDLA"
$ws.Range("H18").Value = "2x NVDLA"

# --- Strip stray leading/trailing whitespace from various text cells ---
$ws.Range("R2").Value = "USB Type-C connector: 2x USB 3.2 Gen2 USB Type-A connector: 2x USB 3.2 Gen2, 2x USB 3.2 Gen1  USB Micro-B connector: USB 2.0"
$ws.Range("U3").Value = "4x UART, 3x SPI, 4x I2S, 8x I2C, 2x CAN, PWM, DMIC & DSPK, GPIOs"
$ws.Range("U4").Value = "4x UART, 3x SPI, 4x I2S, 8x I2C, 2x CAN, PWM, DMIC & DSPK, GPIOs"
$ws.Range("U5").Value = "4x UART, 3x SPI, 4x I2S, 8x I2C, 2x CAN, PWM, DMIC & DSPK, GPIOs"
$ws.Range("N8").Value = "1080p30 supported by 1-2 CPU cores"
$ws.Range("Q8").Value = "M.2 Key M slot with x4 PCIe Gen3 M.2 Key M slot with x2 PCIe Gen3 M.2 Key E slot"
$ws.Range("R8").Value = "USB Type-A Connector: 4x USB 3.2 Gen2  USB Type-C Connector for UFP"
$ws.Range("S8").Value = "1xGbE Connector"
$ws.Range("W8").Value = "100 mm x 79 mmx 21 mm (Height includes feet, carrier board, module, and thermal solution)"
$ws.Range("N9").Value = "1080p30 supported by 1-2 CPU cores"
$ws.Range("D10").Value = "512-core NVIDIA Ampere architecture GPU with 16 Tensor Cores"
$ws.Range("N10").Value = "1080p30 supported by 1-2 CPU cores"
$ws.Range("D14").Value = "512-core NVIDIA Volta architecture GPU with 64 Tensor Cores"
$ws.Range("D15").Value = "512-core NVIDIA Volta architecture GPU with 64 Tensor Cores"
$ws.Range("E15").Value = "1377 MHz"
$ws.Range("D16").Value = "512-core NVIDIA Volta architecture GPU with 64 Tensor Cores"
$ws.Range("E16").Value = "1377 MHz"
$ws.Range("D17").Value = "384-core NVIDIA Volta™ architecture GPU with 48 Tensor Cores"
$ws.Range("E17").Value = "1100 MHz"
$ws.Range("D18").Value = "384-core NVIDIA Volta™ architecture GPU with 48 Tensor Cores"
$ws.Range("E18").Value = "1100 MHz"

# --- Remove the now-redundant "Camera" / " DL Accelerator " columns (X, Y) entirely ---
$ws.Range("X1:Y23").Clear()
